$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user records to append as rows 22-30
$ids      = @(110021, 110022, 110023, 110024, 110025, 110026, 110027, 110028, 110029)
$uins     = @(7316931025, 9137847236, 8428758532, 9804209494, 7105248214, 9316557128, 8103486949, 9601932866, 9317596765)
$names    = @("Magdalena Weber", "Adrienne Hoffman", "Adrienne Mcgee", "Amare Coleman", "Dawson Ibarra", "Elvis Mcmillan", "Steve George", "Colton Elliott", "Carolyn Rodriguez")
$emails   = @("magdalena.weber@xyz.com", "adrienne.hoffman@xyz.com", "adrienne.mcgee@xyz.com", "amare.coleman@xyz.com", "dawson.ibarra@xyz.com", "elvis.mcmillan@xyz.com", "steve.george@xyz.com", "colton.elliott@xyz.com", "carolyn.rodriguez@xyz.com")
$mobiles  = @(932122450, 848488000, 894773246, 956554588, 765455583, 884282274, 971073663, 809908673, 818876429)

# Constant values reused by every data row (identical to existing rows 2-21)
$statusCode = "ACT"
$langCode = "eng"
$lastLoginMethod = "PWD"
$crBy = "superadmin"
$crDtimes = "now()"

$firstRow = 22
$count = $ids.Length

# Populate column-by-column (matching the order the source strings were
# first introduced into the shared-strings table: A, B, C, D, E, F, G, H, I, J, K)
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstRow + $i, 1).Value = $ids[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstRow + $i, 2).Value = $uins[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstRow + $i, 3).Value = $names[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstRow + $i, 4).Value = $emails[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstRow + $i, 5).Value = $mobiles[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstRow + $i, 6).Value = $statusCode
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstRow + $i, 7).Value = $langCode
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstRow + $i, 8).Value = $lastLoginMethod
}
for ($i = 0; $i -lt $count; $i++) {
    $cell = $ws.Cells.Item($firstRow + $i, 9)
    $cell.Value = $true
    $cell.HorizontalAlignment = -4131
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstRow + $i, 10).Value = $crBy
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstRow + $i, 11).Value = $crDtimes
}

# Match the author's final viewport/selection state: rows 22:30 (the newly
# added records) are selected, with the window scrolled so row 16 is on top.
$null = $ws.Range("A22:K30").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

